$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Device (G8): PRN800 -> PX-PR
$ws.Range("G8").Value = "PX-PR"

# Update Panel Type (A8): Pro32xD -> Pro815D
$ws.Range("A8").Value = "Pro815D"

# Update CPU Type (B4): NGC-1609 -> NGC-571/T1396
$ws.Range("B4").Value = "NGC-571/T1396"

# Update 24V PSU load values
$ws.Range("F8").Value = 0.319
$ws.Range("J8").Value = 0.319
$ws.Range("K8").Value = 0.319

# Update the active selection to B8 (was H8)
$ws.Range("B8").Select()
